$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Select()

$ws.Range("C16").Value = 615.32159999999999
$ws.Range("C17").Value = 16490.163199999999
$ws.Range("C18").Value = 32061.004799999999
$ws.Range("C19").Value = 47631.845999999998
$ws.Range("C20").Value = 63202.688000000002

$ws.Range("D16").Value = 609.28
$ws.Range("D17").Value = 16478.080000000002
$ws.Range("D18").Value = 32042.879000000001
$ws.Range("D19").Value = 47607.68
$ws.Range("D20").Value = 63172.480000000003

$ws.Range("D32").Formula = "=C16"
$ws.Range("D33").Formula = "=C17"
$ws.Range("D34").Formula = "=C18"
$ws.Range("D35").Formula = "=C19"
$ws.Range("D36").Formula = "=C20"

$ws.Range("B47").Formula = "=D16"
$ws.Range("B48").Formula = "=D17"
$ws.Range("B49").Formula = "=D18"
$ws.Range("B50").Formula = "=D19"
$ws.Range("B51").Formula = "=D20"

$excel.ActiveWindow.ScrollRow = 17
$ws.Range("B52").Select()
